$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Alternativni tok 1" block (rows 22-31) ---
# Update "Alternativni tok1 1:" row text to describe the full-payment option
$ws.Range("B23").Value = "Korisnik je odabrao opciju plaćanja cjelokupnog iznosa"
$ws.Range("B24").Value = "Na desetom koraku glavnog toka, korisnik je izabrao opciju plaćanja cjelokupnog iznosa"

# The actor switches from "Fizičko lice" to "Pravno lice"
$ws.Range("A27").Value = "Pravno lice"

# Replace the 3-line outcome text with a single new line, then delete the
# now-unused two trailing rows (the block shrinks by 2 rows)
$ws.Range("B28").Value = "1. Nakon utvrđivanja identiteta, pravnom licu se dodjeljuje ključ sale"
$ws.Rows("29:30").Delete()

# --- "Alternativni tok 2" block (rows shifted up by 2: now 30-38) ---
# Update "Alternativni tok1 2:" row text to describe the down-payment option
$ws.Range("B31").Value = "Korisnik je odabrao opciju plaćanja avansa"
$ws.Range("B32").Value = "Na destom koraku glavnog toka, korisnik je izabrao opciju plaćanja avansa"

# The actor switches from "Fizičko lice" to "Pravno lice"
$ws.Range("A35").Value = "Pravno lice"

# Replace the 3-line outcome text with three new lines describing the
# down-payment follow-up flow
$ws.Range("B36").Value = "1. Nakon utvrđivanja identiteta, pravno lice se obavještava o iznosu koji treba doplatiti"
$ws.Range("A37").Value = "2. Doplaćuje iznos usluge"
$ws.Range("B37").Value = ""
$ws.Range("B38").Value = "3. Pravnom licu se dodjeljuje ključ sale"

# Restore the view state (zoom level + selected cell) to match the saved file
$excel.ActiveWindow.Zoom = 80
$ws.Range("B28").Select() | Out-Null
